#
# Weekly CompStat refresh: new crime data collected.
# - Bump the report volume/number and the covered week dates.
# - Shift a week's worth of counts/percentages across the precinct
#   crime-complaints table (rows 14-33) and the historical-perspective
#   labels (rows 39-46).
#

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Force a cell to hold a *text* value (not an auto-converted number),
    # then restore the donor cell's look (number format/alignment/etc.)
    # so the style index matches what a genuine text cell in this sheet uses.
    param($Sheet, $TargetAddr, $DonorAddr, $Text)
    $Sheet.Range($TargetAddr).NumberFormat = "@"
    $Sheet.Range($TargetAddr).Value = $Text
    $Sheet.Range($DonorAddr).Copy()
    $Sheet.Range($TargetAddr).PasteSpecial(-4122)
}

function Set-NumericValue {
    # Force a cell to hold a *numeric* value and borrow an existing
    # numeric cell's formatting so the style matches.
    param($Sheet, $TargetAddr, $DonorAddr, $Number)
    $Sheet.Range($TargetAddr).NumberFormat = "General"
    $Sheet.Range($TargetAddr).Value = $Number
    $Sheet.Range($DonorAddr).Copy()
    $Sheet.Range($TargetAddr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 1. Header: volume number and the covered week's date range moved forward.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# ---------------------------------------------------------------------------
# 2. Cells whose underlying type flips between numeric and text
#    ("***.*" is the sheet's placeholder for an undefined % change, and
#    some counts dropped to/from zero which this report renders as the
#    text "0" rather than a number in a few spots).
# ---------------------------------------------------------------------------
Set-TextValue    $ws "D14" "C14" "0"
Set-TextValue    $ws "E14" "C14" "***.*"
Set-NumericValue $ws "C22" "F22" 1
Set-TextValue    $ws "C29" "C14" "0"
Set-TextValue    $ws "D29" "C14" "0"
Set-TextValue    $ws "E29" "C14" "***.*"
Set-TextValue    $ws "C30" "C14" "0"
Set-TextValue    $ws "D30" "C14" "0"
Set-TextValue    $ws "E30" "C14" "***.*"

# ---------------------------------------------------------------------------
# 3. Plain numeric updates across the crime-complaints table.
# ---------------------------------------------------------------------------
$numericChanges = @(
    @{ Addr = "L14"; Val = -33.333333333333 },
    @{ Addr = "M14"; Val = -66.666666666666 },
    @{ Addr = "D15"; Val = 1 },
    @{ Addr = "F15"; Val = 2 },
    @{ Addr = "H15"; Val = -50 },
    @{ Addr = "J15"; Val = 18 },
    @{ Addr = "K15"; Val = -16.666666666666 },
    @{ Addr = "C16"; Val = 8 },
    @{ Addr = "D16"; Val = 3 },
    @{ Addr = "E16"; Val = 166.666666666667 },
    @{ Addr = "G16"; Val = 13 },
    @{ Addr = "H16"; Val = 38.461538461538 },
    @{ Addr = "I16"; Val = 117 },
    @{ Addr = "J16"; Val = 119 },
    @{ Addr = "K16"; Val = -1.680672268907 },
    @{ Addr = "L16"; Val = -4.098360655737 },
    @{ Addr = "M16"; Val = -6.4 },
    @{ Addr = "N16"; Val = -65.588235294117 },
    @{ Addr = "C17"; Val = 8 },
    @{ Addr = "D17"; Val = 12 },
    @{ Addr = "E17"; Val = -33.333333333333 },
    @{ Addr = "F17"; Val = 59 },
    @{ Addr = "H17"; Val = 78.787878787878 },
    @{ Addr = "I17"; Val = 245 },
    @{ Addr = "J17"; Val = 209 },
    @{ Addr = "K17"; Val = 17.224880382775 },
    @{ Addr = "L17"; Val = 25.641025641025 },
    @{ Addr = "M17"; Val = 114.912280701754 },
    @{ Addr = "N17"; Val = 71.328671328671 },
    @{ Addr = "F18"; Val = 14 },
    @{ Addr = "G18"; Val = 11 },
    @{ Addr = "H18"; Val = 27.272727272727 },
    @{ Addr = "I18"; Val = 117 },
    @{ Addr = "J18"; Val = 86 },
    @{ Addr = "K18"; Val = 36.046511627907 },
    @{ Addr = "L18"; Val = 9.345794392523 },
    @{ Addr = "M18"; Val = -31.578947368421 },
    @{ Addr = "N18"; Val = -83.839779005524 },
    @{ Addr = "C19"; Val = 9 },
    @{ Addr = "D19"; Val = 11 },
    @{ Addr = "E19"; Val = -18.181818181818 },
    @{ Addr = "F19"; Val = 50 },
    @{ Addr = "H19"; Val = -12.280701754386 },
    @{ Addr = "I19"; Val = 359 },
    @{ Addr = "J19"; Val = 420 },
    @{ Addr = "K19"; Val = -14.523809523809 },
    @{ Addr = "L19"; Val = 19.666666666666 },
    @{ Addr = "M19"; Val = 99.444444444444 },
    @{ Addr = "N19"; Val = 31.021897810219 },
    @{ Addr = "C20"; Val = 8 },
    @{ Addr = "E20"; Val = -38.461538461538 },
    @{ Addr = "F20"; Val = 38 },
    @{ Addr = "G20"; Val = 42 },
    @{ Addr = "H20"; Val = -9.523809523809 },
    @{ Addr = "I20"; Val = 239 },
    @{ Addr = "J20"; Val = 218 },
    @{ Addr = "K20"; Val = 9.633027522935 },
    @{ Addr = "L20"; Val = -5.15873015873 },
    @{ Addr = "M20"; Val = 109.649122807018 },
    @{ Addr = "N20"; Val = -73.024830699774 },
    @{ Addr = "C21"; Val = 36 },
    @{ Addr = "D21"; Val = 42 },
    @{ Addr = "E21"; Val = -14.285714285714 },
    @{ Addr = "F21"; Val = 181 },
    @{ Addr = "G21"; Val = 161 },
    @{ Addr = "H21"; Val = 12.422360248447 },
    @{ Addr = "I21"; Val = 1094 },
    @{ Addr = "J21"; Val = 1073 },
    @{ Addr = "K21"; Val = 1.957129543336 },
    @{ Addr = "L21"; Val = 10.393541876892 },
    @{ Addr = "M21"; Val = 51.733703190013 },
    @{ Addr = "N21"; Val = -54.149203688181 },
    @{ Addr = "F22"; Val = 2 },
    @{ Addr = "I22"; Val = 7 },
    @{ Addr = "K22"; Val = -22.222222222222 },
    @{ Addr = "L22"; Val = -22.222222222222 },
    @{ Addr = "M22"; Val = -41.666666666666 },
    @{ Addr = "C23"; Val = 1 },
    @{ Addr = "E23"; Val = -50 },
    @{ Addr = "F23"; Val = 13 },
    @{ Addr = "G23"; Val = 6 },
    @{ Addr = "H23"; Val = 116.666666666667 },
    @{ Addr = "I23"; Val = 55 },
    @{ Addr = "J23"; Val = 55 },
    @{ Addr = "K23"; Val = 0 },
    @{ Addr = "L23"; Val = -17.910447761194 },
    @{ Addr = "M23"; Val = 96.428571428571 },
    @{ Addr = "C24"; Val = 29 },
    @{ Addr = "D24"; Val = 26 },
    @{ Addr = "E24"; Val = 11.538461538461 },
    @{ Addr = "F24"; Val = 177 },
    @{ Addr = "G24"; Val = 87 },
    @{ Addr = "H24"; Val = 103.448275862069 },
    @{ Addr = "I24"; Val = 789 },
    @{ Addr = "J24"; Val = 673 },
    @{ Addr = "K24"; Val = 17.236255572065 },
    @{ Addr = "L24"; Val = 6.910569105691 },
    @{ Addr = "M24"; Val = 96.268656716417 },
    @{ Addr = "C25"; Val = 12 },
    @{ Addr = "D25"; Val = 8 },
    @{ Addr = "E25"; Val = 50 },
    @{ Addr = "F25"; Val = 73 },
    @{ Addr = "G25"; Val = 29 },
    @{ Addr = "H25"; Val = 151.724137931034 },
    @{ Addr = "I25"; Val = 302 },
    @{ Addr = "J25"; Val = 273 },
    @{ Addr = "K25"; Val = 10.62271062271 },
    @{ Addr = "L25"; Val = 1.342281879194 },
    @{ Addr = "C26"; Val = 6 },
    @{ Addr = "D26"; Val = 18 },
    @{ Addr = "E26"; Val = -66.666666666666 },
    @{ Addr = "F26"; Val = 62 },
    @{ Addr = "G26"; Val = 61 },
    @{ Addr = "H26"; Val = 1.639344262295 },
    @{ Addr = "I26"; Val = 341 },
    @{ Addr = "J26"; Val = 272 },
    @{ Addr = "K26"; Val = 25.367647058823 },
    @{ Addr = "L26"; Val = 29.657794676806 },
    @{ Addr = "M26"; Val = 14.046822742474 },
    @{ Addr = "D27"; Val = 1 },
    @{ Addr = "F27"; Val = 2 },
    @{ Addr = "G27"; Val = 4 },
    @{ Addr = "H27"; Val = -50 },
    @{ Addr = "J27"; Val = 21 },
    @{ Addr = "K27"; Val = -14.285714285714 },
    @{ Addr = "L27"; Val = 0 },
    @{ Addr = "C28"; Val = 1 },
    @{ Addr = "E28"; Val = -50 },
    @{ Addr = "F28"; Val = 7 },
    @{ Addr = "G28"; Val = 4 },
    @{ Addr = "H28"; Val = 75 },
    @{ Addr = "I28"; Val = 34 },
    @{ Addr = "J28"; Val = 30 },
    @{ Addr = "K28"; Val = 13.333333333333 },
    @{ Addr = "L28"; Val = -12.820512820512 },
    @{ Addr = "L29"; Val = -60 },
    @{ Addr = "M29"; Val = -66.666666666666 },
    @{ Addr = "N29"; Val = -81.818181818181 },
    @{ Addr = "L30"; Val = -66.666666666666 },
    @{ Addr = "M30"; Val = -70 },
    @{ Addr = "N30"; Val = -84.210526315789 },
    @{ Addr = "L31"; Val = -50 }
)

foreach ($change in $numericChanges) {
    $ws.Range($change.Addr).Value = $change.Val
}

Write-Host "Applied $($numericChanges.Count) numeric updates plus header/type-change edits."
